$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values
$ws.Range("B3").Value = "0.2.2"
$ws.Range("B8").Value = "2024-09-11T16:17:59-05:00"
$ws.Range("B10").Value = "MITRE (https://www.mitre.org)"

# Insert a new row for "Jurisdiction" before the Description row (old row 11)
$ws.Rows.Item(11).Insert()

$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""
